$wb = $excel.ActiveWorkbook

# Work on Sheet2: delete the erroneous sorted row (row 6), which shifts
# rows 7:14 up to become rows 6:13.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows.Item(6).Delete()

# Select the block that now occupies A6:C13, anchored at A6, and make
# Sheet2 the active sheet/tab.
$ws2.Activate()
$ws2.Range("A6:C13").Select()

# Sheet1 should no longer be the selected tab; Sheet2's selection above
# already makes Sheet2 the active (and only selected) sheet.
